$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B values (rows 1-20) to match re-layout of the table
$values = @(0,1,1,0,0,0,0,0,0,0,0,0,1,0,0,0,0,0,0,0)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Re-layout: selection moved to D13
$ws.Range("D13").Select()
